$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $newValue
    $c.ClearFormats()
}

Set-TextValue 'D2' '29.902.27'
Set-TextValue 'E2' '  +0.71%  '
Set-TextValue 'D3' '1.627.65'
Set-TextValue 'E3' '  +1.48%  '
Set-TextValue 'E4' '  +0.05%  '
Set-TextValue 'D5' '214.23'
Set-TextValue 'E5' '  +0.93%  '
Set-TextValue 'E6' '  +0.67%  '
Set-TextValue 'E7' '  +0.06%  '
Set-TextValue 'E8' '  +7.91%  '
Set-TextValue 'E9' '  +2.94%  '
Set-TextValue 'D10' '0.0612'
Set-TextValue 'E10' '  +1.92%  '
Set-TextValue 'E11' '  +0.89%  '
Set-TextValue 'D12' '1.860.94'
Set-TextValue 'E12' '  +1.47%  '
Set-TextValue 'D13' '1.630.66'
Set-TextValue 'E13' '  +1.98%  '
Set-TextValue 'D14' '0.574'
Set-TextValue 'E14' '  +6.88%  '
Set-TextValue 'E15' '  +4.90%  '
Set-TextValue 'D16' '29.964.61'
Set-TextValue 'E16' '  +0.98%  '
Set-TextValue 'D17' '8.92'
Set-TextValue 'E17' '  +17.54%  '
Set-TextValue 'D18' '64.74'
Set-TextValue 'D19' '242.98'
Set-TextValue 'E19' '  +0.62%  '
Set-TextValue 'E20' '  +2.06%  '
Set-TextValue 'E22' '  +3.74%  '
Set-TextValue 'D23' '9.65'
Set-TextValue 'E23' '  +4.51%  '
Set-TextValue 'E24' '  +2.12%  '
Set-TextValue 'D25' '157.33'
Set-TextValue 'E25' '  +1.55%  '
Set-TextValue 'D26' '15.72'
Set-TextValue 'E26' '  +2.44%  '
Set-TextValue 'E27' '  +2.74%  '
Set-TextValue 'E28' '  +3.22%  '
Set-TextValue 'E29' '  +0.05%  '
Set-TextValue 'E30' '  +3.11%  '
Set-TextValue 'E31' '  +5.45%  '
Set-TextValue 'E33' '  +3.15%  '
Set-TextValue 'D34' '1.423.59'
Set-TextValue 'E34' '  -0.53%  '
Set-TextValue 'E35' '  +6.91%  '
Set-TextValue 'D36' '1.04'
Set-TextValue 'E36' '  +0.36%  '
Set-TextValue 'E37' '  +1.97%  '
Set-TextValue 'E38' '  -0.33%  '
Set-TextValue 'E39' '  +3.50%  '
Set-TextValue 'D40' '0.560'
Set-TextValue 'E40' '  +4.01%  '
Set-TextValue 'E41' '  +0.86%  '
Set-TextValue 'E42' '  +3.87%  '
Set-TextValue 'D43' '0.0499'
Set-TextValue 'E43' '  +2.27%  '
Set-TextValue 'D44' '54.35'
Set-TextValue 'E44' '  -0.53%  '
Set-TextValue 'D45' '69.49'
Set-TextValue 'E45' '  +5.54%  '
Set-TextValue 'E46' '  +7.88%  '
Set-TextValue 'D47' '0.999'
Set-TextValue 'E47' '  +0.07%  '
Set-TextValue 'D48' '5.40'
Set-TextValue 'D49' '1.768.43'
Set-TextValue 'E49' '  +1.34%  '
Set-TextValue 'D50' '88.94'
Set-TextValue 'E50' '  +2.86%  '
Set-TextValue 'E51' '  +5.57%  '
